$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 10853.4
$ws.Range("I34").Value = 9066.75
$ws.Range("J34").Value = 18000
$ws.Range("K34").Value = 9066.75
$ws.Range("L34").Value = 18000
$ws.Range("M34").Value = -8863.75
$ws.Range("N34").Value = -18406
$ws.Range("H36").Value = 10853.4
$ws.Range("I36").Value = 9066.75
$ws.Range("J36").Value = 18000
$ws.Range("K36").Value = 9066.75
$ws.Range("L36").Value = 18000
$ws.Range("M36").Value = -8351.75
$ws.Range("N36").Value = -19430
$ws.Range("H40").Value = 166668160
$ws.Range("I40").Value = 2250.5
$ws.Range("K40").Value = 2250.5
$ws.Range("M40").Value = -2075.5
$ws.Range("H41").Value = 2687.5
$ws.Range("I41").Value = 3341.1
$ws.Range("J41").Value = 1598.1666
$ws.Range("K41").Value = 3341.1
$ws.Range("L41").Value = 1598.1666
$ws.Range("M41").Value = -2901.1
$ws.Range("N41").Value = -2478.1666
$ws.Range("H43").Value = 5862.636
$ws.Range("J43").Value = 5449.5
$ws.Range("L43").Value = 5449.5
$ws.Range("N43").Value = -5587.5
$ws.Range("H76").Value = 8651.75
$ws.Range("I76").Value = 10252.333
$ws.Range("J76").Value = 3850
$ws.Range("K76").Value = 10252.333
$ws.Range("L76").Value = 3850
$ws.Range("M76").Value = -9937.333000000001
$ws.Range("N76").Value = -4480
$ws.Range("H79").Value = 8651.75
$ws.Range("I79").Value = 10252.333
$ws.Range("J79").Value = 3850
$ws.Range("K79").Value = 10252.333
$ws.Range("L79").Value = 3850
$ws.Range("M79").Value = -9160.333000000001
$ws.Range("N79").Value = -6034
$ws.Range("H116").Value = 12833.571
$ws.Range("I116").Value = 6700.5713
$ws.Range("J116").Value = 15900.071
$ws.Range("K116").Value = 6700.5713
$ws.Range("L116").Value = 15900.071
$ws.Range("M116").Value = -3258.5713
$ws.Range("N116").Value = -22784.071
$ws.Range("H135").Value = 1385.8572
$ws.Range("I135").Value = 648.4828
$ws.Range("J135").Value = 4949.8335
$ws.Range("K135").Value = 5836.3452
$ws.Range("L135").Value = 44548.5015
$ws.Range("M135").Value = -3301.3452
$ws.Range("N135").Value = -49618.5015
$ws.Range("H138").Value = 5374.836
$ws.Range("J138").Value = 6329.8184
$ws.Range("L138").Value = 18989.4552
$ws.Range("N138").Value = -29269.4552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7876.256
$ws.Range("I32").Value = 7701.8887
$ws.Range("J32").Value = 22000
$ws.Range("K32").Value = 7701.8887
$ws.Range("L32").Value = 22000
$ws.Range("M32").Value = -7414.8887
$ws.Range("N32").Value = -22574
$ws.Range("H97").Value = 1417.3928
$ws.Range("J97").Value = 2356.1
$ws.Range("L97").Value = 2356.1
$ws.Range("N97").Value = -3348.1
$ws.Range("H132").Value = 2883.0625
$ws.Range("I132").Value = 2947.2068
$ws.Range("J132").Value = 2263
$ws.Range("K132").Value = 8841.6204
$ws.Range("L132").Value = 6789
$ws.Range("M132").Value = -6311.6204
$ws.Range("N132").Value = -11849
$ws.Range("H137").Value = 99949.17999999999
$ws.Range("J137").Value = 100944.1
$ws.Range("L137").Value = 100944.1
$ws.Range("N137").Value = -111144.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5550.4873
$ws.Range("J20").Value = 3598.1177
$ws.Range("L20").Value = 3598.1177
$ws.Range("N20").Value = -4092.1177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 741.8946999999999
$ws.Range("I16").Value = 593.1177
$ws.Range("K16").Value = 593.1177
$ws.Range("M16").Value = -306.1177
$ws.Range("H31").Value = 25644298
$ws.Range("I31").Value = 30305880
$ws.Range("K31").Value = 30305880
$ws.Range("M31").Value = -30305585
$ws.Range("H32").Value = 8504.5
$ws.Range("I32").Value = 8504.5
$ws.Range("K32").Value = 8504.5
$ws.Range("M32").Value = -8188.5
$ws.Range("H34").Value = 25644298
$ws.Range("I34").Value = 30305880
$ws.Range("K34").Value = 30305880
$ws.Range("M34").Value = -30305678
$ws.Range("H99").Value = 19275.285
$ws.Range("I99").Value = 11857.571
$ws.Range("J99").Value = 22984.143
$ws.Range("K99").Value = 11857.571
$ws.Range("L99").Value = 22984.143
$ws.Range("M99").Value = -10359.571
$ws.Range("N99").Value = -25980.143
$ws.Range("H103").Value = 26274.54
$ws.Range("I103").Value = 12196.375
$ws.Range("K103").Value = 12196.375
$ws.Range("M103").Value = -11024.375
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180
$ws.Range("H113").Value = 741.8946999999999
$ws.Range("I113").Value = 593.1177
$ws.Range("K113").Value = 593.1177
$ws.Range("M113").Value = 1576.8823
$ws.Range("H126").Value = 19275.285
$ws.Range("I126").Value = 11857.571
$ws.Range("J126").Value = 22984.143
$ws.Range("K126").Value = 35572.713
$ws.Range("L126").Value = 68952.429
$ws.Range("M126").Value = -33102.713
$ws.Range("N126").Value = -73892.429
$ws.Range("H141").Value = 358751.28
$ws.Range("I141").Value = 109765.336
$ws.Range("J141").Value = 452121
$ws.Range("K141").Value = 109765.336
$ws.Range("L141").Value = 452121
$ws.Range("M141").Value = -104585.336
$ws.Range("N141").Value = -462481

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4650
$ws.Range("I34").Value = 6000
$ws.Range("J34").Value = 4380
$ws.Range("K34").Value = 18000
$ws.Range("L34").Value = 13140
$ws.Range("M34").Value = -17916
$ws.Range("N34").Value = -13308
$ws.Range("H122").Value = 55899.5
$ws.Range("I122").Value = 165500
$ws.Range("J122").Value = 1099.25
$ws.Range("K122").Value = 1489500
$ws.Range("L122").Value = 9893.25
$ws.Range("M122").Value = -1487050
$ws.Range("N122").Value = -14793.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 82629.71000000001
$ws.Range("I59").Value = 86899
$ws.Range("J59").Value = 80922
$ws.Range("K59").Value = 86899
$ws.Range("L59").Value = 80922
$ws.Range("M59").Value = -86316
$ws.Range("N59").Value = -82088
$ws.Range("H80").Value = 3835.3333
$ws.Range("I80").Value = 2930.2856
$ws.Range("J80").Value = 7003
$ws.Range("K80").Value = 2930.2856
$ws.Range("L80").Value = 7003
$ws.Range("M80").Value = -1932.2856
$ws.Range("N80").Value = -8999
$ws.Range("H83").Value = 3835.3333
$ws.Range("I83").Value = 2930.2856
$ws.Range("J83").Value = 7003
$ws.Range("K83").Value = 14651.428
$ws.Range("L83").Value = 35015
$ws.Range("M83").Value = -9659.428
$ws.Range("N83").Value = -44999
$ws.Range("H126").Value = 13589996
$ws.Range("I126").Value = 19182984
$ws.Range("K126").Value = 57548952
$ws.Range("M126").Value = -57546482

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 26405140
$ws.Range("I22").Value = 44005430
$ws.Range("K22").Value = 44005430
$ws.Range("M22").Value = -44005135
$ws.Range("H27").Value = 26405140
$ws.Range("I27").Value = 44005430
$ws.Range("K27").Value = 44005430
$ws.Range("M27").Value = -44005323
$ws.Range("H40").Value = 5192.973
$ws.Range("I40").Value = 4510.2
$ws.Range("K40").Value = 4510.2
$ws.Range("M40").Value = -4374.2
$ws.Range("H43").Value = 8000
$ws.Range("I43").Value = 8000
$ws.Range("K43").Value = 8000
$ws.Range("M43").Value = -7807
$ws.Range("H61").Value = 3188.3809
$ws.Range("I61").Value = 3077.8
$ws.Range("K61").Value = 3077.8
$ws.Range("M61").Value = -2875.8
$ws.Range("H113").Value = 3188.3809
$ws.Range("I113").Value = 3077.8
$ws.Range("K113").Value = 3077.8
$ws.Range("M113").Value = -907.8000000000002
$ws.Range("H122").Value = 3175.5605
$ws.Range("I122").Value = 3162.6785
$ws.Range("K122").Value = 9488.0355
$ws.Range("M122").Value = -7038.0355
$ws.Range("H132").Value = 4252.273
$ws.Range("I132").Value = 3150.6428
$ws.Range("K132").Value = 9451.928400000001
$ws.Range("M132").Value = -6921.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3727.4092
$ws.Range("I136").Value = 3896
$ws.Range("J136").Value = 3366.1428
$ws.Range("K136").Value = 11688
$ws.Range("L136").Value = 10098.4284
$ws.Range("M136").Value = -9138
$ws.Range("N136").Value = -15198.4284
